$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.002.47"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.338.09"
$ws.Range("E3").Value = "  +4.19%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'309.45"
$ws.Range("E5").Value = "  +3.81%  "
$ws.Range("D6").Value = "'107.63"
$ws.Range("E6").Value = "  -4.99%  "
$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D9").Value = "'0.623"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "'43.25"
$ws.Range("E10").Value = "  -5.29%  "
$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'8.92"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "'1.07"
$ws.Range("E13").Value = "  +18.15%  "
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "'16.27"
$ws.Range("E15").Value = "  +6.42%  "
$ws.Range("D16").Value = "2.691.23"
$ws.Range("E16").Value = "  +4.33%  "
$ws.Range("D17").Value = "2.421.13"
$ws.Range("E17").Value = "  +7.28%  "
$ws.Range("D18").Value = "42.914.16"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'7.21"
$ws.Range("E20").Value = "  -4.44%  "
$ws.Range("D21").Value = "'75.41"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").Value = "'2.52"
$ws.Range("E23").Value = "  +7.58%  "
$ws.Range("D24").Value = "'249.64"
$ws.Range("E24").Value = "  +7.68%  "
$ws.Range("D25").Value = "'8.94"
$ws.Range("E25").Value = "  -5.62%  "
$ws.Range("D26").Value = "'11.87"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D29").Value = "'38.60"
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("D30").Value = "'22.40"
$ws.Range("E30").Value = "  +4.94%  "
$ws.Range("D31").Value = "'173.94"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'3.16"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("D33").Value = "'0.0907"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "'5.82"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").Value = "'5.01"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").Value = "'0.130"
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").Value = "'4.10"
$ws.Range("E38").Value = "  -5.63%  "
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("E40").Value = "  +7.48%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'71.89"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.48"
$ws.Range("E42").Value = "  +10.88%  "
$ws.Range("D43").Value = "'0.233"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "'12.36"
$ws.Range("E45").Value = "  -7.23%  "
$ws.Range("D46").Value = "'5.67"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "'9.21"
$ws.Range("E47").Value = "  +5.78%  "
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").Value = "'0.0996"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").Value = "'70.08"
$ws.Range("E51").Value = "  +2.77%  "

# Reset style on cells that required a leading apostrophe to stay text-typed,
# so no stray number-format style is left attached to the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
